# Commit: "Moved to readxl, fixed column headings."
#
# The "Game N" / "Good Session" column headers in row 1 of Sheet1 get
# dotted names (as R's readxl / make.names would produce), e.g.
# "Game 1" -> "Game.1", "Good Session" -> "Good.Session".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G1").Value = "Game.1"
$ws.Range("H1").Value = "Game.2"
$ws.Range("I1").Value = "Game.3"
$ws.Range("J1").Value = "Game.4"
$ws.Range("K1").Value = "Game.5"
$ws.Range("L1").Value = "Game.6"
$ws.Range("M1").Value = "Good.Session"

# Selection moved from the data area (H28) back up to the header row.
$ws.Range("A1:N1").Select()
